$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.286.57'
$ws.Range("E2").Value = '  +4.16%  '
$ws.Range("D3").Value = '3.485.52'
$ws.Range("E3").Value = '  +3.68%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.59'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.64'
$ws.Range("E6").Value = '  +7.44%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +1.60%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  +4.80%  '
$ws.Range("E11").Value = '  +4.75%  '
$ws.Range("D12").Value = '4.078.38'
$ws.Range("E12").Value = '  +3.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.81'
$ws.Range("E13").Value = '  +7.40%  '
$ws.Range("D15").Value = '3.486.56'
$ws.Range("E15").Value = '  +3.66%  '
$ws.Range("E16").Value = '  +4.35%  '
$ws.Range("D17").Value = '63.296.48'
$ws.Range("E17").Value = '  +3.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.28'
$ws.Range("E18").Value = '  +3.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.38'
$ws.Range("E19").Value = '  +6.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.42'
$ws.Range("E20").Value = '  +6.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '392.69'
$ws.Range("E21").Value = '  +2.47%  '
$ws.Range("E22").Value = '  +3.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.24'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("E25").Value = '  +9.55%  '
$ws.Range("D26").Value = '3.628.08'
$ws.Range("E26").Value = '  +3.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.185'
$ws.Range("E27").Value = '  -2.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.90'
$ws.Range("E28").Value = '  +11.05%  '
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.26'
$ws.Range("E30").Value = '  +5.75%  '
$ws.Range("E31").Value = '  +2.52%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.83'
$ws.Range("E34").Value = '  +3.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '32.71'
$ws.Range("E35").Value = '  +30.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.19'
$ws.Range("E36").Value = '  +5.63%  '
$ws.Range("E37").Value = '  +9.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '171.83'
$ws.Range("E38").Value = '  +2.38%  '
$ws.Range("E39").Value = '  +10.06%  '
$ws.Range("D40").Value = '3.522.94'
$ws.Range("E40").Value = '  +3.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0770'
$ws.Range("E41").Value = '  +2.31%  '
$ws.Range("E42").Value = '  +4.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.75'
$ws.Range("E43").Value = '  +7.98%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.50'
$ws.Range("E44").Value = '  +4.02%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.44'
$ws.Range("E45").Value = '  +0.37%  '
$ws.Range("E46").Value = '  +10.71%  '
$ws.Range("D47").Value = '2.614.75'
$ws.Range("E47").Value = '  +6.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.84'
$ws.Range("E48").Value = '  +7.92%  '
$ws.Range("E49").Value = '  +17.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.77'
$ws.Range("E50").Value = '  +2.56%  '
$ws.Range("E51").Value = '  +5.77%  '
